$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add "ID" header in B1 (inherits style from B1, already styled s="1")
$ws.Range("B1").Value = "ID"

# Fill in the IDs for each department row
$ws.Range("B2").Value = 123123
$ws.Range("B3").Value = 1001
$ws.Range("B4").Value = 2032

# Update the active selection to B4
$ws.Range("B4").Select()
